$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '29.093.58'
$ws.Range('E2').Value = '  -3.56%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.847.92'
$ws.Range('E3').Value = '  -2.60%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9995'
$ws.Range('E4').Value = '  +0.12%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.7057'
$ws.Range('E5').Value = '  -5.36%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '238.17'
$ws.Range('E6').Value = '  -1.95%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.9995'
$ws.Range('E7').Value = '  +0.28%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3051'
$ws.Range('E8').Value = '  -3.83%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07500'
$ws.Range('E9').Value = '  +3.35%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '23.37'
$ws.Range('E10').Value = '  -7.00%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.08132'
$ws.Range('E11').Value = '  -2.79%  '
$ws.Range('B12').Value = 'WrappedEther'
$ws.Range('C12').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.847.33'
$ws.Range('E12').Value = '  -2.77%  '
$ws.Range('B13').Value = 'Polygon'
$ws.Range('C13').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.7248'
$ws.Range('E13').Value = '  -5.01%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.229'
$ws.Range('E14').Value = '  -3.99%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '89.21'
$ws.Range('E15').Value = '  -4.22%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '29.126.31'
$ws.Range('E16').Value = '  -3.38%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '5.794'
$ws.Range('E17').Value = '  -6.21%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '239.82'
$ws.Range('E18').Value = '  -4.18%  '
$ws.Range('E19').Value = '  -2.59%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '13.07'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.001'
$ws.Range('E21').Value = '  +0.40%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '2.093.13'
$ws.Range('E22').Value = '  -1.79%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.9998'
$ws.Range('E23').Value = '  +0.20%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '7.556'
$ws.Range('E24').Value = '  -6.18%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.1465'
$ws.Range('E25').Value = '  -7.69%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '8.974'
$ws.Range('E26').Value = '  -3.77%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '160.97'
$ws.Range('E27').Value = '  -1.98%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '18.03'
$ws.Range('E28').Value = '  -4.05%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.940'
$ws.Range('E29').Value = '  -5.97%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.383'
$ws.Range('E30').Value = '  -6.17%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.568'
$ws.Range('E31').Value = '  -0.98%  '
$ws.Range('E32').Value = '  -2.81%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.009'
$ws.Range('E33').Value = '  -5.64%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.05170'
$ws.Range('E34').Value = '  -4.17%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.187'
$ws.Range('E35').Value = '  -5.63%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.031'
$ws.Range('E36').Value = '  +3.12%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.7062'
$ws.Range('E37').Value = '  -7.64%  '
$ws.Range('E38').Value = '  -2.71%  '
$ws.Range('E39').Value = '  -5.67%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.679'
$ws.Range('E40').Value = '  -3.38%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.9095'
$ws.Range('E41').Value = '  +4.24%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.986'
$ws.Range('E42').Value = '  -1.56%  '
$ws.Range('B43').Value = 'Maker'
$ws.Range('C43').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.078.47'
$ws.Range('E43').Value = '  -2.09%  '
$ws.Range('B44').Value = 'TheSandbox'
$ws.Range('C44').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.4301'
$ws.Range('E44').Value = '  -6.03%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '70.07'
$ws.Range('E45').Value = '  -3.92%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.9992'
$ws.Range('E46').Value = '  +0.14%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '102.29'
$ws.Range('E47').Value = '  -2.28%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.752'
$ws.Range('E48').Value = '  -6.56%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '7.068'
$ws.Range('E49').Value = '  -7.44%  '
$ws.Range('B50').Value = 'RocketPoolETH'
$ws.Range('C50').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.982.80'
$ws.Range('E50').Value = '  -3.29%  '
$ws.Range('B51').Value = 'EnergySwap'
$ws.Range('C51').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '9.176'
$ws.Range('E51').Value = '  -4.82%  '
